$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.92"
$ws.Range("E2").Value = "'0.69%"
$ws.Range("D3").Value = "'41.27"
$ws.Range("E3").Value = "'4.91%"
$ws.Range("D4").Value = "'5.123"
$ws.Range("E4").Value = "'0.62%"
$ws.Range("D5").Value = "'0.07645"
$ws.Range("E5").Value = "'-0.63%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.622"
$ws.Range("E6").Value = "'-0.19%"
$ws.Range("B7").Value = "BTSEToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D7").Value = "'2.455"
$ws.Range("E7").Value = "'0.90%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9061"
$ws.Range("E8").Value = "'-0.87%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1124"
$ws.Range("E9").Value = "'10.65%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1809"
$ws.Range("E10").Value = "'3.04%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09082"
$ws.Range("E11").Value = "'-2.08%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04241"
$ws.Range("E12").Value = "'-4.41%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1052"
$ws.Range("E13").Value = "'-0.51%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001251"
$ws.Range("E14").Value = "'-0.57%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005728"
$ws.Range("E15").Value = "'-1.68%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.346"
$ws.Range("E16").Value = "'-0.49%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.268"
$ws.Range("E17").Value = "'0.52%"
$ws.Range("E18").Value = "'0.42%"
$ws.Range("D19").Value = "'6.693"
$ws.Range("E19").Value = "'-4.78%"
$ws.Range("D20").Value = "'0.1359"
$ws.Range("E20").Value = "'0.86%"
$ws.Range("D22").Value = "'0.04051"
$ws.Range("E22").Value = "'-2.21%"
$ws.Range("D23").Value = "'0.001265"
$ws.Range("E23").Value = "'5.52%"
$ws.Range("D24").Value = "'0.004041"
$ws.Range("D25").Value = "'0.0001271"
$ws.Range("E25").Value = "'-2.28%"
$ws.Range("D26").Value = "'0.0003744"
$ws.Range("D38").Value = "'0.02434"
$ws.Range("E38").Value = "'0.19%"
$ws.Range("D39").Value = "'0.05261"
$ws.Range("E39").Value = "'1.58%"
$ws.Range("D40").Value = "'0.007795"
$ws.Range("E40").Value = "'-1.38%"
$ws.Range("D41").Value = "'0.1302"
$ws.Range("E41").Value = "'-1.25%"
$ws.Range("D42").Value = "'0.006527"
$ws.Range("E42").Value = "'-8.97%"
$ws.Range("D43").Value = "'0.001951"
$ws.Range("E43").Value = "'0.07%"
$ws.Range("D44").Value = "'0.007568"
$ws.Range("E44").Value = "'-9.59%"
$ws.Range("D45").Value = "'0.3085"
$ws.Range("E45").Value = "'0.93%"
$ws.Range("D46").Value = "'0.00006778"
$ws.Range("E46").Value = "'5.66%"
$ws.Range("E47").Value = "'0.01%"
$ws.Range("D48").Value = "'0.06988"
$ws.Range("E48").Value = "'1,478.06%"
$ws.Range("E49").Value = "'39.94%"
$ws.Range("E50").Value = "'0.01%"
$ws.Range("E51").Value = "'0.01%"
